{"js": "// The document text content got \"rotated\" across 8 paragraphs: the body\n// text that used to live in one slot (Objetivos / Docente / Programa\n// resumido / Programa / M\u00e9todo-value / Crit\u00e9rio-value / Norma-value /\n// Bibliografia) now lives in a different slot, while every paragraph's\n// style/position stays exactly where it was (verified: same 16\n// paragraphs, same styles, same order, before and after).\n//\n// We therefore address each paragraph by its stable index and overwrite\n// its text with the new target content (a literal derived straight from\n// the diff), rather than trying to generically \"move\" ranges around.\n// \"\\u000b\" (vertical tab) is Office.js's in-memory marker for a soft\n// line break (<w:br/>), so joining segments with it reproduces the\n// original multi-line-inside-one-run structure.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- Paragraph 5: \"Objetivos\" body -------------------------------------\nparagraphs.items[5].insertText(\n  \"Sistemas de Informa\u00e7\u00e3o. Projeto de Sistemas de Informa\u00e7\u00e3o. Tecnologia de Informa\u00e7\u00e3o. Processo de Desenvolvimento de Sistema de Informa\u00e7\u00e3o.\",\n  \"Replace\"\n);\n\n// --- Paragraph 7: \"Docente(s) Respons\u00e1vel(eis)\" list item body ---------\nparagraphs.items[7].insertText(\n  \"Oferecer ao aluno uma vis\u00e3o geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementa\u00e7\u00e3o de Sistemas de Informa\u00e7\u00e3o no sentido de capacit\u00e1-lo analisa e projetar tais sistemas\",\n  \"Replace\"\n);\n\n// --- Paragraph 9: \"Programa resumido\" body ------------------------------\nconst programaSegments = [\n  \"1. Sistemas de Informa\u00e7\u00e3o\",\n  \"1.1. Sistemas de Processamento de Informa\u00e7\u00f5es;\",\n  \"1.2. Sistemas de Informa\u00e7\u00f5es Gerenciais;\",\n  \"1.3. Sistema de Apoio \u00e0 Decis\u00e3o;\",\n  \"1.4. Sistemas de Informa\u00e7\u00e3o no Com\u00e9rcio Eletr\u00f4nico;\",\n  \"1.5. Sistemas de Informa\u00e7\u00e3o em Cadeia de Suprimentos;\",\n  \"1.6. Sistemas inteligentes nos neg\u00f3cios;\",\n  \"1.7. Sistemas estrat\u00e9gicos. \",\n  \"2. Projeto de Sistemas de Informa\u00e7\u00e3o.\",\n  \"2.1. Especifica\u00e7\u00e3o das Sa\u00eddas;\",\n  \"2.2. Especifica\u00e7\u00e3o dos Arquivos;\",\n  \"2.3. Especifica\u00e7\u00e3o das Entradas;\",\n  \"2.4. Especifica\u00e7\u00e3o do Processamento.\",\n  \"3. Tecnologia de Informa\u00e7\u00e3o.\",\n  \"3.1. Evolu\u00e7\u00e3o da Computa\u00e7\u00e3o;\",\n  \"3.2. Recursos Computacionais.\",\n  \"4. Processo de Desenvolvimento de Sistemas de Informa\u00e7\u00e3o.\",\n  \"4.1. Defini\u00e7\u00e3o do Neg\u00f3cio;\",\n  \"4.2. Identifica\u00e7\u00e3o do Problema e/ou Oportunidades;\",\n  \"4.3. Sele\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;\",\n  \"4.4. Implementa\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;\",\n  \"4.5. Avalia\u00e7\u00e3o da Efic\u00e1cia do Sistema de Informa\u00e7\u00e3o;\",\n];\nparagraphs.items[9].insertText(programaSegments.join(\"\\u000b\"), \"Replace\");\n\n// --- Paragraph 11: \"Programa\" body --------------------------------------\nparagraphs.items[11].insertText(\n  \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios.\",\n  \"Replace\"\n);\n\n// --- Paragraph 13: \"Avalia\u00e7\u00e3o\" list item (M\u00e9todo/Crit\u00e9rio/Norma) -------\n// This paragraph keeps its bold \"M\u00e9todo: \" / \"Crit\u00e9rio: \" / \"Norma de\n// recupera\u00e7\u00e3o: \" label runs untouched; only the plain-text value run\n// that follows each label changes. Each value is unique text within\n// this paragraph, so we can scope the search to the paragraph itself\n// and swap each value run in place with insertText(..., \"Replace\").\n// We go in reverse (Norma -> Crit\u00e9rio -> M\u00e9todo) so that the text\n// being searched for is never something we *just* inserted earlier in\n// this same paragraph (the three old values are a simple chain: the\n// new M\u00e9todo value is the old Crit\u00e9rio value, the new Crit\u00e9rio value is\n// the old Norma value, and the new Norma value is the former\n// Bibliografia text).\nconst p13 = paragraphs.items[13];\n\nconst bibliografiaSegments = [\n  \"HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004\",\n  \"LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gest\u00e3o integrada de processos e da tecnologia da informa\u00e7\u00e3o. S\u00e3o Paulo:Atlas, 2006.\",\n  \"LAURINDO, F.J.B. Tecnologia da Informa\u00e7\u00e3o: Efic\u00e1cia nas Organiza\u00e7\u00f5es. S\u00e3o Paulo, Editora Futura, 2002.\",\n  \"STAIR, R.M., Princ\u00edpios de Sistema de Informa\u00e7\u00e3o: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.\",\n  \"TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.\",\n  \"TURBAN, E., RAIANER JR, K., POTTER, R. E., Administra\u00e7\u00e3o de Tecnologia da Informa\u00e7\u00e3o: Teoria e Pr\u00e1tica\\u201d, S\u00e3o Paulo, Editora Campus, 2003.\",\n];\n\nasync function replaceValueInParagraph(paragraph, oldText, newText) {\n  const found = paragraph.search(oldText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for \" +\n        JSON.stringify(oldText) +\n        \" but found \" +\n        found.items.length\n    );\n  }\n  found.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Norma de recupera\u00e7\u00e3o value: old \"M\u00e9dia aritm\u00e9tica da nota final...\" -> bibliography text\nawait replaceValueInParagraph(\n  p13,\n  \"M\u00e9dia aritm\u00e9tica da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recupera\u00e7\u00e3o.\",\n  bibliografiaSegments.join(\"\\u000b\")\n);\n\n// Crit\u00e9rio value: old \"M\u00e9dia Aritm\u00e9tica das atividades...\" -> \"M\u00e9dia aritm\u00e9tica da nota final...\"\nawait replaceValueInParagraph(\n  p13,\n  \"M\u00e9dia Aritm\u00e9tica das atividades avaliativas realizadas.\",\n  \"M\u00e9dia aritm\u00e9tica da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recupera\u00e7\u00e3o.\"\n);\n\n// M\u00e9todo value: old \"Aulas expositivas...\" -> \"M\u00e9dia Aritm\u00e9tica das atividades...\"\nawait replaceValueInParagraph(\n  p13,\n  \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios.\",\n  \"M\u00e9dia Aritm\u00e9tica das atividades avaliativas realizadas.\"\n);\n\n// --- Paragraph 15: \"Bibliografia\" body ----------------------------------\nparagraphs.items[15].insertText(\"5840917 - Fabricio Maciel Gomes\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The document text content got \"rotated\" across 8 paragraphs: the body\n# text that used to live in one slot (Objetivos / Docente / Programa\n# resumido / Programa / Metodo-value / Criterio-value / Norma-value /\n# Bibliografia) now lives in a different slot, while every paragraph's\n# style/position stays exactly where it was (16 paragraphs, same styles,\n# same order, before and after).\n#\n# We therefore address each paragraph by its stable 1-based Paragraphs\n# index and overwrite its text with the new target content (a literal\n# derived straight from the diff), rather than trying to generically\n# \"move\" ranges around. [char]11 is Word's manual-line-break character\n# (produces <w:br/> when written into a Range), so joining segments with\n# it reproduces the original multi-line-inside-one-run structure.\n\n$d = $word.ActiveDocument\n$lineBreak = [char]11\n\n# --- Paragraph 6: \"Objetivos\" body --------------------------------------\n$d.Paragraphs.Item(6).Range.Text = \"Sistemas de Informa\u00e7\u00e3o. Projeto de Sistemas de Informa\u00e7\u00e3o. Tecnologia de Informa\u00e7\u00e3o. Processo de Desenvolvimento de Sistema de Informa\u00e7\u00e3o.\"\n\n# --- Paragraph 8: \"Docente(s) Respons\u00e1vel(eis)\" list item body ---------\n$d.Paragraphs.Item(8).Range.Text = \"Oferecer ao aluno uma vis\u00e3o geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementa\u00e7\u00e3o de Sistemas de Informa\u00e7\u00e3o no sentido de capacit\u00e1-lo analisa e projetar tais sistemas\"\n\n# --- Paragraph 10: \"Programa resumido\" body -----------------------------\n$programaSegments = @(\n    \"1. Sistemas de Informa\u00e7\u00e3o\",\n    \"1.1. Sistemas de Processamento de Informa\u00e7\u00f5es;\",\n    \"1.2. Sistemas de Informa\u00e7\u00f5es Gerenciais;\",\n    \"1.3. Sistema de Apoio \u00e0 Decis\u00e3o;\",\n    \"1.4. Sistemas de Informa\u00e7\u00e3o no Com\u00e9rcio Eletr\u00f4nico;\",\n    \"1.5. Sistemas de Informa\u00e7\u00e3o em Cadeia de Suprimentos;\",\n    \"1.6. Sistemas inteligentes nos neg\u00f3cios;\",\n    \"1.7. Sistemas estrat\u00e9gicos. \",\n    \"2. Projeto de Sistemas de Informa\u00e7\u00e3o.\",\n    \"2.1. Especifica\u00e7\u00e3o das Sa\u00eddas;\",\n    \"2.2. Especifica\u00e7\u00e3o dos Arquivos;\",\n    \"2.3. Especifica\u00e7\u00e3o das Entradas;\",\n    \"2.4. Especifica\u00e7\u00e3o do Processamento.\",\n    \"3. Tecnologia de Informa\u00e7\u00e3o.\",\n    \"3.1. Evolu\u00e7\u00e3o da Computa\u00e7\u00e3o;\",\n    \"3.2. Recursos Computacionais.\",\n    \"4. Processo de Desenvolvimento de Sistemas de Informa\u00e7\u00e3o.\",\n    \"4.1. Defini\u00e7\u00e3o do Neg\u00f3cio;\",\n    \"4.2. Identifica\u00e7\u00e3o do Problema e/ou Oportunidades;\",\n    \"4.3. Sele\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;\",\n    \"4.4. Implementa\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;\",\n    \"4.5. Avalia\u00e7\u00e3o da Efic\u00e1cia do Sistema de Informa\u00e7\u00e3o;\"\n)\n$d.Paragraphs.Item(10).Range.Text = [string]::Join($lineBreak, $programaSegments)\n\n# --- Paragraph 12: \"Programa\" body --------------------------------------\n$d.Paragraphs.Item(12).Range.Text = \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios.\"\n\n# --- Paragraph 14: \"Avalia\u00e7\u00e3o\" list item (M\u00e9todo/Crit\u00e9rio/Norma) -------\n# This paragraph keeps its bold \"M\u00e9todo: \" / \"Crit\u00e9rio: \" / \"Norma de\n# recupera\u00e7\u00e3o: \" label runs untouched; only the plain-text value run\n# that follows each label changes. Each value is unique text within\n# this paragraph, so we scope Find/Replace to that paragraph's Range\n# (wdReplaceOne) and swap each value run in place, leaving every\n# sibling run (labels, bold formatting, line breaks) alone.\n#\n# We go in reverse (Norma -> Crit\u00e9rio -> M\u00e9todo) so the text being\n# searched for is never something we *just* inserted earlier in this\n# same paragraph: the old values chain together (new M\u00e9todo value is\n# the old Crit\u00e9rio value; new Crit\u00e9rio value is the old Norma value;\n# new Norma value is the former Bibliografia text), so replacing from\n# the bottom up keeps every search hit unique at the moment it runs.\n\n$bibliografiaSegments = @(\n    \"HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004\",\n    \"LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gest\u00e3o integrada de processos e da tecnologia da informa\u00e7\u00e3o. S\u00e3o Paulo:Atlas, 2006.\",\n    \"LAURINDO, F.J.B. Tecnologia da Informa\u00e7\u00e3o: Efic\u00e1cia nas Organiza\u00e7\u00f5es. S\u00e3o Paulo, Editora Futura, 2002.\",\n    \"STAIR, R.M., Princ\u00edpios de Sistema de Informa\u00e7\u00e3o: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.\",\n    \"TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.\",\n    \"TURBAN, E., RAIANER JR, K., POTTER, R. E., Administra\u00e7\u00e3o de Tecnologia da Informa\u00e7\u00e3o: Teoria e Pr\u00e1tica\u201d, S\u00e3o Paulo, Editora Campus, 2003.\"\n)\n$bibliografiaJoined = [string]::Join($lineBreak, $bibliografiaSegments)\n\n# NOTE: a successful Find/Replace collapses the Range down to just the\n# replaced text, so we re-fetch a fresh Paragraphs.Item(14).Range before\n# each call instead of reusing one Range across all three replacements.\n\n# Norma de recupera\u00e7\u00e3o value: old \"M\u00e9dia aritm\u00e9tica da nota final...\" -> bibliography text\n$d.Paragraphs.Item(14).Range.Find.Execute(\"M\u00e9dia aritm\u00e9tica da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recupera\u00e7\u00e3o.\", $false, $false, $false, $false, $false, $true, 1, $false, $bibliografiaJoined, 1)\n\n# Crit\u00e9rio value: old \"M\u00e9dia Aritm\u00e9tica das atividades...\" -> \"M\u00e9dia aritm\u00e9tica da nota final...\"\n$d.Paragraphs.Item(14).Range.Find.Execute(\"M\u00e9dia Aritm\u00e9tica das atividades avaliativas realizadas.\", $false, $false, $false, $false, $false, $true, 1, $false, \"M\u00e9dia aritm\u00e9tica da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recupera\u00e7\u00e3o.\", 1)\n\n# M\u00e9todo value: old \"Aulas expositivas...\" -> \"M\u00e9dia Aritm\u00e9tica das atividades...\"\n$d.Paragraphs.Item(14).Range.Find.Execute(\"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios.\", $false, $false, $false, $false, $false, $true, 1, $false, \"M\u00e9dia Aritm\u00e9tica das atividades avaliativas realizadas.\", 1)\n\n# --- Paragraph 16: \"Bibliografia\" body -----------------------------------\n$d.Paragraphs.Item(16).Range.Text = \"5840917 - Fabricio Maciel Gomes\"\n"}
